$d = $word.ActiveDocument

# New numbered use-case items to append after "View Spending Category Report"
$items = @(
    "View Income Source Report",
    "View Cash Flow Report",
    "View Account Listing Report",
    "Create User Profile"
)

foreach ($item in $items) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $rng = $lastPara.Range
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.Text = $item
}

# Trailing empty paragraph: still styled as ListParagraph, but not part of the numbered list
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $lastPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>'
$newPara.Range.InsertXML($xmlFrag)
